# Daily attendance processing - 2025-10-18 04:43:52
# Reorders the "Recorded By" (column G) values in the
# "Session Analysis Results" sheet so that real user / email entries are
# listed before generic "System"/"system" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "backup@backdoor.com, system, System"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    7   = "admin@admin.com, System"
    10  = "dnasr281@gmail.com, System"
    18  = "dnasr281@gmail.com, System"
    19  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, system, System"
    31  = "backup@backdoor.com, System"
    32  = "backup@backdoor.com, System"
    34  = "admin@admin.com, System"
    37  = "dnasr281@gmail.com, System"
    45  = "dnasr281@gmail.com, System"
    46  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, system, System"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    61  = "admin@admin.com, System"
    64  = "dnasr281@gmail.com, System"
    72  = "dnasr281@gmail.com, System"
    73  = "dnasr281@gmail.com, System"
    83  = "backup@backdoor.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    86  = "dnasr281@gmail.com, System"
    90  = "dnasr281@gmail.com, admin@admin.com"
    97  = "dnasr281@gmail.com, System"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    112 = "dnasr281@gmail.com, System"
    116 = "dnasr281@gmail.com, admin@admin.com"
    123 = "dnasr281@gmail.com, System"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    138 = "dnasr281@gmail.com, System"
    142 = "dnasr281@gmail.com, admin@admin.com"
    149 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
